$d = $word.ActiveDocument

# The document currently ends with an empty ListParagraph-styled paragraph
# (numId 1, ilvl 0) right before the sectPr. We fill it with the first new
# question, then append four more list paragraphs with the same formatting.

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertAfter("Should shirts be shipped on separate line items by PO or just combine and apply to the oldest PO first?")
$r.Collapse(0)
$r.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertAfter("When knowing what cut order each shipment came from, is it important to know specific to the box when there is more than 1 cut order represented?")
$r.Collapse(0)
$r.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertAfter("How do you handle quality at each stage?")
$r.Collapse(0)
$r.InsertParagraphAfter()

# This paragraph needs a spell-check proofErr split around "shirtstyleSKUs",
# so build it from a literal OOXML fragment instead of plain InsertAfter.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">Figure out way to handle color of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>shirtstyleSKUs</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
'@
$r.InsertXML($xml)

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.InsertAfter("Are certain SKUs only allowed in certain colors?  How are colors limited?  By customer?  Is there just a separate SKU in the catalog for every color combo of every size of every style?")
